# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the Granada data block (row 109),
# pushing the existing rows 109:130 down to 110:131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 109 (shifts 109:130 -> 110:131)
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new record
$ws.Cells.Item(109, 1).Value  = 10
$ws.Cells.Item(109, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(109, 3).Value  = "La Araucanía"
$ws.Cells.Item(109, 4).Value  = 44736
$ws.Cells.Item(109, 5).Value  = 9
$ws.Cells.Item(109, 6).Value  = "Fruta"
$ws.Cells.Item(109, 7).Value  = 100104
$ws.Cells.Item(109, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(109, 9).Value  = 100104001
$ws.Cells.Item(109, 10).Value = "Granada"
$ws.Cells.Item(109, 11).Value = "Wonderfull"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 40
$ws.Cells.Item(109, 14).Value = 15000
$ws.Cells.Item(109, 15).Value = 15000
$ws.Cells.Item(109, 16).Value = 15000
$ws.Cells.Item(109, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(109, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(109, 19).Value = 1000
$ws.Cells.Item(109, 20).Value = 15
